$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("B2").Value = 8
$ws.Range("D2").Value = 11
$ws.Range("C18").Value = 16
$ws.Range("D18").Value = 3

# Update the active selection (was E1:F18 with active cell E1, now just C2)
$ws.Range("C2").Select()
